$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4376.971
$ws.Range("I64").Value = 3578.2
$ws.Range("J64").Value = 4510.1
$ws.Range("K64").Value = 3578.2
$ws.Range("L64").Value = 4510.1
$ws.Range("M64").Value = -3330.2
$ws.Range("N64").Value = -5006.1

$ws.Range("H67").Value = 4376.971
$ws.Range("I67").Value = 3578.2
$ws.Range("J67").Value = 4510.1
$ws.Range("K67").Value = 3578.2
$ws.Range("L67").Value = 4510.1
$ws.Range("M67").Value = -2720.2
$ws.Range("N67").Value = -6226.1

$ws.Range("H76").Value = 4300
$ws.Range("I76").Value = 4266.6665
$ws.Range("K76").Value = 4266.6665
$ws.Range("M76").Value = -3951.6665

$ws.Range("H79").Value = 4300
$ws.Range("I79").Value = 4266.6665
$ws.Range("K79").Value = 4266.6665
$ws.Range("M79").Value = -3174.6665

$ws.Range("H129").Value = 940.0714
$ws.Range("J129").Value = 1093.8372
$ws.Range("L129").Value = 3281.5116
$ws.Range("N129").Value = -13281.5116

$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -60120

$ws.Range("H136").Value = 60650
$ws.Range("J136").Value = 60650
$ws.Range("L136").Value = 60650
$ws.Range("N136").Value = -70850

$ws.Range("H139").Value = 80625
$ws.Range("J139").Value = 80625
$ws.Range("L139").Value = 80625
$ws.Range("N139").Value = -90905

$ws.Range("H140").Value = 84166
$ws.Range("J140").Value = 83123.336
$ws.Range("L140").Value = 83123.336
$ws.Range("N140").Value = -93483.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13875.041
$ws.Range("I32").Value = 13801.016
$ws.Range("K32").Value = 13801.016
$ws.Range("M32").Value = -13514.016

$ws.Range("H63").Value = 10108.091
$ws.Range("I63").Value = 11190
$ws.Range("J63").Value = 9206.5
$ws.Range("K63").Value = 11190
$ws.Range("L63").Value = 9206.5
$ws.Range("M63").Value = -10504
$ws.Range("N63").Value = -10578.5

$ws.Range("H66").Value = 10108.091
$ws.Range("I66").Value = 11190
$ws.Range("J66").Value = 9206.5
$ws.Range("K66").Value = 55950
$ws.Range("L66").Value = 46032.5
$ws.Range("M66").Value = -52518
$ws.Range("N66").Value = -52896.5

$ws.Range("H88").Value = 3115.8572
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 3135.1667
$ws.Range("K88").Value = 3000
$ws.Range("L88").Value = 3135.1667
$ws.Range("M88").Value = -2594
$ws.Range("N88").Value = -3947.1667

$ws.Range("H91").Value = 3115.8572
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 3135.1667
$ws.Range("K91").Value = 3000
$ws.Range("L91").Value = 3135.1667
$ws.Range("M91").Value = -1596
$ws.Range("N91").Value = -5943.1667

$ws.Range("H112").Value = 22523.125
$ws.Range("J112").Value = 22523.125
$ws.Range("L112").Value = 22523.125
$ws.Range("N112").Value = -25477.125

$ws.Range("H133").Value = 43241.668
$ws.Range("J133").Value = 43241.668
$ws.Range("L133").Value = 43241.668
$ws.Range("N133").Value = -48301.668

$ws.Range("H134").Value = 51758
$ws.Range("J134").Value = 51758
$ws.Range("L134").Value = 51758
$ws.Range("N134").Value = -61898

$ws.Range("H139").Value = 84746
$ws.Range("J139").Value = 84746
$ws.Range("L139").Value = 84746
$ws.Range("N139").Value = -95026

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3618.7856
$ws.Range("I105").Value = 3169.3635
$ws.Range("K105").Value = 3169.3635
$ws.Range("M105").Value = -1422.3635

$ws.Range("H107").Value = 42642.848
$ws.Range("I107").Value = 66675.875
$ws.Range("J107").Value = 4190
$ws.Range("K107").Value = 66675.875
$ws.Range("L107").Value = 4190
$ws.Range("M107").Value = -64755.875
$ws.Range("N107").Value = -8030

$ws.Range("H132").Value = 72764.82000000001
$ws.Range("J132").Value = 72764.82000000001
$ws.Range("L132").Value = 72764.82000000001
$ws.Range("N132").Value = -82884.82000000001

$ws.Range("H135").Value = 74476.664
$ws.Range("J135").Value = 74476.664
$ws.Range("L135").Value = 74476.664
$ws.Range("N135").Value = -84616.664

$ws.Range("H137").Value = 45378.57
$ws.Range("J137").Value = 45378.57
$ws.Range("L137").Value = 45378.57
$ws.Range("N137").Value = -55578.57

$ws.Range("H138").Value = 42961.25
$ws.Range("J138").Value = 42961.25
$ws.Range("L138").Value = 42961.25
$ws.Range("N138").Value = -53241.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 76854
$ws.Range("J135").Value = 133708
$ws.Range("L135").Value = 133708
$ws.Range("N135").Value = -143848

$ws.Range("H140").Value = 74576
$ws.Range("J140").Value = 74576
$ws.Range("L140").Value = 74576
$ws.Range("N140").Value = -84936

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 2175
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2175
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 6525
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -7897

$ws.Range("H65").Value = 2175
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2175
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 19575
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -26439

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6305.45
$ws.Range("I70").Value = 5450
$ws.Range("J70").Value = 6519.3125
$ws.Range("K70").Value = 5450
$ws.Range("L70").Value = 6519.3125
$ws.Range("M70").Value = -5180
$ws.Range("N70").Value = -7059.3125

$ws.Range("H73").Value = 6305.45
$ws.Range("I73").Value = 5450
$ws.Range("J73").Value = 6519.3125
$ws.Range("K73").Value = 5450
$ws.Range("L73").Value = 6519.3125
$ws.Range("M73").Value = -4514
$ws.Range("N73").Value = -8391.3125

$ws.Range("H80").Value = 3499.3076
$ws.Range("I80").Value = 2857.8572
$ws.Range("J80").Value = 4247.6665
$ws.Range("K80").Value = 2857.8572
$ws.Range("L80").Value = 4247.6665
$ws.Range("M80").Value = -1859.8572
$ws.Range("N80").Value = -6243.6665

$ws.Range("H83").Value = 3499.3076
$ws.Range("I83").Value = 2857.8572
$ws.Range("J83").Value = 4247.6665
$ws.Range("K83").Value = 14289.286
$ws.Range("L83").Value = 21238.3325
$ws.Range("M83").Value = -9297.286
$ws.Range("N83").Value = -31222.3325

$ws.Range("H111").Value = 29900
$ws.Range("J111").Value = 29900
$ws.Range("L111").Value = 29900
$ws.Range("N111").Value = -36034

$ws.Range("H135").Value = 71419.2
$ws.Range("J135").Value = 71419.2
$ws.Range("L135").Value = 71419.2
$ws.Range("N135").Value = -81559.2

$ws.Range("H138").Value = 66400
$ws.Range("J138").Value = 66400
$ws.Range("L138").Value = 66400
$ws.Range("N138").Value = -76680

$ws.Range("H140").Value = 44410
$ws.Range("J140").Value = 44410
$ws.Range("L140").Value = 44410
$ws.Range("N140").Value = -54770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 161290.75
$ws.Range("J69").Value = 161290.75
$ws.Range("L69").Value = 161290.75
$ws.Range("N69").Value = -162912.75

$ws.Range("H72").Value = 161290.75
$ws.Range("J72").Value = 161290.75
$ws.Range("L72").Value = 483872.25
$ws.Range("N72").Value = -491984.25

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H133").Value = 58798.523
$ws.Range("J133").Value = 58798.523
$ws.Range("L133").Value = 58798.523
$ws.Range("N133").Value = -63858.523

$ws.Range("H135").Value = 165498.17
$ws.Range("J135").Value = 165498.17
$ws.Range("L135").Value = 165498.17
$ws.Range("N135").Value = -175638.17

$ws.Range("H137").Value = 87000
$ws.Range("J137").Value = 87000
$ws.Range("L137").Value = 87000
$ws.Range("N137").Value = -97200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 54175.715
$ws.Range("J135").Value = 54175.715
$ws.Range("L135").Value = 54175.715
$ws.Range("N135").Value = -64315.715

$ws.Range("H137").Value = 39329.668
$ws.Range("J137").Value = 39329.668
$ws.Range("L137").Value = 39329.668
$ws.Range("N137").Value = -49529.668

$ws.Range("H139").Value = 59423.332
$ws.Range("J139").Value = 59423.332
$ws.Range("L139").Value = 59423.332
$ws.Range("N139").Value = -69703.33199999999
